$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures that are stored as plain text in the
# source workbook (not real numbers - e.g. '27.260.94', '5.360').
# Force text format on every Price cell we touch so Excel doesn't
# silently reinterpret the new value as a numeric literal.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "27.266.29"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.908.50"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "307.37"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.5264"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "0.07294"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "21.96"
$ws.Range("E10").Value = "  +3.80%  "
$ws.Range("D11").Value = "0.9023"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "0.08168"
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("D13").Value = "96.00"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "5.362"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "1.452.33"
$ws.Range("E15").Value = "  -23.83%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "0.000008668"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "14.77"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "27.300.52"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "5.109"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "10.82"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "6.520"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "149.95"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").Value = "2.304"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "18.25"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "1.742"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "116.62"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").Value = "4.854"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "4.853"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "0.09258"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "0.8298"
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("D33").Value = "0.05071"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "1.229"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("D35").Value = "2.998"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").Value = "3.350"
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("E37").Value = "  +4.93%  "
$ws.Range("D38").Value = "0.5797"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "1.075"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "9.221"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").Value = "6.585"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "116.64"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "0.1525"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").Value = "10.21"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "38.91"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").Value = "0.06210"
$ws.Range("E50").Value = "  +4.14%  "
$ws.Range("E51").Value = "  +0.77%  "
